$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("model_14_5_0", 0.5010782282294184, 0.3791495606328492, 0.3917534915799327, 0.4477410905295015, 60.17011260986328, 94.49362182617188, 36.31206130981445, 67.11407470703125),
    @("model_14_5_1", 0.5122750834436902, 0.3583494960628629, 0.4237981594061286, 0.4413578063043367, 58.81977081298828, 97.65940093994141, 34.39900588989258, 67.88980865478516),
    @("model_14_5_2", 0.5317061992022282, 0.4122970122918594, 0.3899932546234153, 0.4693121976685155, 56.47637176513672, 89.44857025146484, 36.41714477539062, 64.49261474609375),
    @("model_14_5_3", 0.5404194849343491, 0.3501516123105428, 0.4718104776777189, 0.4470216793759235, 55.42554473876953, 98.90712738037109, 31.53269195556641, 67.20149993896484),
    @("model_14_5_4", 0.5410780975414631, 0.3005365506773902, 0.4661037619969374, 0.4128056336692332, 55.34611129760742, 106.4585494995117, 31.87338066101074, 71.35964965820312),
    @("model_14_5_5", 0.5717831119780701, 0.3574807712641276, 0.5401807992196365, 0.4676866370233626, 51.64308547973633, 97.7916259765625, 27.45101165771484, 64.69015502929688),
    @("model_14_5_11", 0.5723071582572985, 0.3078130582154748, 0.1544761545927709, 0.3455896828925454, 51.57988357543945, 105.3510589599609, 50.4774169921875, 79.52816772460938),
    @("model_14_5_9", 0.5823706731643306, 0.3331708866462276, 0.3880416617807921, 0.4163974235730167, 50.36621856689453, 101.4915924072266, 36.53365325927734, 70.92315673828125),
    @("model_14_5_8", 0.5881247841330518, 0.3410415169206469, 0.4766023918389489, 0.442089009402578, 49.67227554321289, 100.2936859130859, 31.24661636352539, 67.80094146728516),
    @("model_14_5_7", 0.5886548588049501, 0.3393866515743771, 0.4938428585284789, 0.4449773289099012, 49.60834503173828, 100.5455551147461, 30.21736907958984, 67.44994354248047),
    @("model_14_5_6", 0.5910898876767523, 0.3456349737453623, 0.5592085369645083, 0.464231152874453, 49.31467819213867, 99.59455871582031, 26.31506538391113, 65.11008453369141),
    @("model_14_5_18", 0.593398236751963, 0.2580415519884932, 0.4504085688605136, 0.381001555309985, 49.03629302978516, 112.9263076782227, 32.81037521362305, 75.22468566894531),
    @("model_14_5_12", 0.5942031478790619, 0.2990034513669816, 0.430647170996982, 0.4035924742312659, 48.9392204284668, 106.6918869018555, 33.99011993408203, 72.47930145263672),
    @("model_14_5_16", 0.5945671373100367, 0.2706326728773504, 0.4462712844753584, 0.388393505802792, 48.89532470703125, 111.0099258422852, 33.057373046875, 74.32637023925781),
    @("model_14_5_17", 0.5960926295059297, 0.2723927272675806, 0.4542532680141328, 0.3914057259748337, 48.71135330200195, 110.7420425415039, 32.58084869384766, 73.96030426025391),
    @("model_14_5_13", 0.5962891782073554, 0.3057080839644775, 0.4356218581596593, 0.4091878988010755, 48.68764495849609, 105.6714401245117, 33.69313430786133, 71.79930114746094),
    @("model_14_5_15", 0.5975831235476825, 0.2893974432227321, 0.4495798754731949, 0.4016001204619458, 48.53159332275391, 108.153923034668, 32.85984802246094, 72.72142028808594),
    @("model_14_5_10", 0.597704652195519, 0.3366048039421951, 0.4451673428492053, 0.4318802956374783, 48.51693725585938, 100.9689559936523, 33.12327575683594, 69.04157257080078),
    @("model_14_5_14", 0.5997646965575385, 0.3116277826452079, 0.4522606734635393, 0.4169594158510636, 48.26849365234375, 104.7704696655273, 32.69980621337891, 70.8548583984375),
    @("model_14_5_19", 0.6002080052318077, 0.2788187579834575, 0.4762574281715917, 0.4007532291321236, 48.21503067016602, 109.7640075683594, 31.26720809936523, 72.8243408203125)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowData[$j]
    }
}

$ws.Range("A22:I24").Delete()
